# Basic card layout for the game workbook:
#  - rename "Deck" -> "Opportunity"
#  - add a new blank "Alien" sheet after it
#  - lay out the Opportunity card table (headers + 3 sample rows)

$wb = $excel.ActiveWorkbook

# --- sheets -----------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Opportunity"

$alien = $wb.Worksheets.Add($null, $ws)
$alien.Name = "Alien"

# the old 3x3 "Deck" grid leaves a stray value (old DEF column) that the
# new layout doesn't reuse - clear it before laying out the new table
$ws.Range("C3").ClearContents()

# --- header row (row 1) ------------------------------------------------
$headers = @{
    "A1" = "Name";
    "B1" = "Tags";
    "C1" = "Consume1";
    "D1" = "Consume1Icon";
    "E1" = "Consume2";
    "F1" = "Consume2Icon";
    "G1" = "Required1";
    "H1" = "Required1Icon";
    "I1" = "Required2";
    "J1" = "Required2Icon";
    "K1" = "Action1";
    "L1" = "Action2";
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# bold header cells (A,B,K,L keep plain bold; C:J are bold + centered)
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("K1:L1").Font.Bold = $true
$ws.Range("C1:J1").Font.Bold = $true
$ws.Range("C1:J1").HorizontalAlignment = -4108

# --- row 2: Sieze Opportunity -------------------------------------------
$ws.Range("A2").Value = "Sieze Opportunity"
$ws.Range("B2").Value = "Action"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "electric"
$ws.Range("K2").Value = "Take an opportunity card from your pile and place it on your tableau."
$ws.Range("C2:D2").HorizontalAlignment = -4108

# --- row 3: Expand Options ----------------------------------------------
$ws.Range("A3").Value = "Expand Options"
$ws.Range("B3").Value = "Buff"
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = "electric"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = "glass-heart"
$ws.Range("K3").Value = "Deal +1 Opportunities"
$ws.Range("G3:J3").HorizontalAlignment = -4108

# --- row 4: Basic Factories ----------------------------------------------
$ws.Range("A4").Value = "Basic Factories"
$ws.Range("B4").Value = "Action, Manufacturing"

# --- column widths (best effort match) -----------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.28515625
$ws.Columns.Item(2).ColumnWidth = 20.85546875
$ws.Columns.Item(3).ColumnWidth = 10.28515625
$ws.Columns.Item(4).ColumnWidth = 15.42578125
$ws.Columns.Item(5).ColumnWidth = 10.28515625
$ws.Columns.Item(6).ColumnWidth = 14.140625
$ws.Columns.Item(7).ColumnWidth = 14.140625
$ws.Columns.Item(8).ColumnWidth = 14.140625
$ws.Columns.Item(9).ColumnWidth = 14.140625
$ws.Columns.Item(10).ColumnWidth = 14.140625
$ws.Columns.Item(11).ColumnWidth = 63.140625

# --- selection / active sheet --------------------------------------------
[void]$ws.Range("D4").Select()
$ws.Activate()
